# DOMA-6936: add "Верифицирован" (Is verified) column to the contacts import example sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Create the new column H by copying the formatting from column G ---
#        (header style, data-row style and the blank bordered template rows 7-50)
$ws.Range("G1:G50").Copy()
$ws.Range("H1:H50").PasteSpecial(-4122)  # xlPasteFormats

# Match the F:G column width for the new column H
$ws.Columns("H").ColumnWidth = $ws.Columns("G").ColumnWidth

# --- 2. Fill in the header and the values for the new "Верифицирован" column ---
$ws.Range("H1").Value = "Верифицирован"
$ws.Range("H2").Value = "Да"
$ws.Range("H3").Value = "Нет"
# H4 stays empty - no verification value for that sample contact
$ws.Range("H5").Value = "да"
$ws.Range("H6").Value = "нет"

# --- 3. Fix the typo in the example e-mail on row 5 (test@example.com -> ttest@example.com) ---
$ws.Range("F5").Value = "ttest@example.com"

# Rebuild the hyperlinks; F5 keeps the same mailto: target as before, only its
# displayed text is corrected to "ttest@example.com" (matches the source edit)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:test@example.com", "", "", "test1@example.com")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:test@example.com", "", "", "test2@example.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:test@example.com", "", "", "ttest@example.com")

# Adding hyperlinks re-applies Excel's builtin "Hyperlink" look (blue/underline) to F2:F5;
# restore the original plain cell formatting (style copied from an untouched column-G cell)
$ws.Range("G2").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)  # xlPasteFormats
